$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.489.03'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.111.34'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '334.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5258'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4497'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.51'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +12.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09056'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.41%  '
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.46'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.098.62'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.793'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.828'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.81'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001130'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06624'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.45'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.325'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.540.51'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.40'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.356'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.348.96'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.84%  '
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '163.43'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '132.91'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.200'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1075'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.667'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.169'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.933'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.94%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.59'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +11.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02583'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06840'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.583'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.60%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2313'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.65%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.81'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6943'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.248'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.382'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.002'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '14.12'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6408'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.670'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.250'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.56%  '
$ws.Range('E50').Value = '  +2.95%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '83.51'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.38%  '
